# Iteration og faseplan - add iteration plan for E1
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1 ("Overordnet projektplan"): add "møde med hans" activities for row 6-8 ---
$ws1.Range("H6").Value = "forberedelse til møde med hans"
$ws1.Range("H7").Value = "møde med hans"
$ws1.Range("I6").Value = "20 min"
$ws1.Range("H8").Value = "implementation af feedback fra hans"
$ws1.Range("I7").Value = "20 min"
$ws1.Range("I8").Value = "20 min"

# --- sheet2 ("Milepæle"): insert a Plan section above the existing Milepæl table ---
# move the existing milestone rows (3-7) down to rows 15-19
$ws2.Range("A15").Value = $ws2.Range("A3").Value()
$ws2.Range("B15").Value = $ws2.Range("B3").Value()
$ws2.Range("B16").Value = $ws2.Range("B4").Value()
$ws2.Range("B17").Value = $ws2.Range("B5").Value()
$ws2.Range("B18").Value = $ws2.Range("B6").Value()
$ws2.Range("B19").Value = $ws2.Range("B7").Value()
$ws2.Range("B15:B19").Font.Italic = $true
$ws2.Range("B15:B19").HorizontalAlignment = -4131
$ws2.Range("B15:B19").IndentLevel = 1

# clear the old rows 3:7 content
$ws2.Range("A3:F7").Clear()

# write the new Plan section
$ws2.Range("A3").Value = "Plan"
$ws2.Range("C3").Value = "Lave ikke funktionelle krav (FURPS)"
$ws2.Range("C4").Value = "Identificer alle use cases."
$ws2.Range("C5").Value = "Formel beskriv UC3 – opret lånetilbud"
$ws2.Range("C6").Value = "Sekvens diagram for UC1 og UC2"
$ws2.Range("C7").Value = "Opdatere klasse diagram"
$ws2.Range("C8").Value = "Implementere UC1 og UC2."

$ws2.Columns("C").ColumnWidth = 35.140625

# --- add new empty sheet "Ark1" at the end ---
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "Ark1"

# --- restore sheet1's own selection (it stays the non-active sheet) ---
$ws1.Select()
$ws1.Range("F26").Select()

# --- sheet2 becomes the active/selected sheet, with its own remembered selection ---
$ws2.Select()
$ws2.Range("C12").Select()
